# "Working on find transaction test scenarios"
#
# 1. Sign up sheet: Username sample value "username1" -> "username113"
# 2. Sign up sheet: cursor/selection moved to I2 (the Username cell)
# 3. General Data sheet: cursor/selection moved to E5
# 4. General Data sheet: new "Transaction date" value added at D2
#    (a date, formatted as a short date) to go with the new
#    "Transaction date from/to" headers already on that sheet.

$wb = $excel.ActiveWorkbook

# --- Sign up sheet -------------------------------------------------
$wsSignup = $wb.Worksheets.Item("Sign up")
$wsSignup.Activate() | Out-Null

$wsSignup.Range("I2").Value = "username113"

$wsSignup.Range("I2").Select() | Out-Null

# --- General Data sheet ---------------------------------------------
$wsGeneral = $wb.Worksheets.Item("General Data")
$wsGeneral.Activate() | Out-Null

$wsGeneral.Range("D2").Value = 45772
$wsGeneral.Range("D2").NumberFormat = "[$-404]b2"

$wsGeneral.Range("E5").Select() | Out-Null

# Leave "General Data" as the active/selected sheet, matching the
# workbook's tabSelected state.
$wsGeneral.Activate() | Out-Null
